# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Faerie Profits workbook
# (currentAveragePrice / profit columns per leve row), per sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 6671199.5
$ws.Range("J40").Value = 5090.909
$ws.Range("L40").Value = 5090.909
$ws.Range("N40").Value = -5440.909
# Row 61
$ws.Range("H61").Value = 221.71428
$ws.Range("I61").Value = 221.71428
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 665.14284
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -493.14284
$ws.Range("N61").Value = $null
# Row 75
$ws.Range("H75").Value = 39457
$ws.Range("J75").Value = 39457
$ws.Range("L75").Value = 39457
$ws.Range("N75").Value = -41329
# Row 78
$ws.Range("H78").Value = 39457
$ws.Range("J78").Value = 39457
$ws.Range("L78").Value = 118371
$ws.Range("N78").Value = -127731
# Row 86
$ws.Range("H86").Value = 2499.5
$ws.Range("I86").Value = 2666.3333
$ws.Range("J86").Value = 1999
$ws.Range("K86").Value = 2666.3333
$ws.Range("L86").Value = 1999
$ws.Range("M86").Value = -1543.3333
$ws.Range("N86").Value = -4245
# Row 89
$ws.Range("H89").Value = 2499.5
$ws.Range("I89").Value = 2666.3333
$ws.Range("J89").Value = 1999
$ws.Range("K89").Value = 13331.6665
$ws.Range("L89").Value = 9995
$ws.Range("M89").Value = -7715.666499999999
$ws.Range("N89").Value = -21227
# Row 94
$ws.Range("H94").Value = 200
$ws.Range("I94").Value = 200
$ws.Range("K94").Value = 200
$ws.Range("M94").Value = 251
# Row 113
$ws.Range("H113").Value = 6899.143
$ws.Range("I113").Value = 3798.5
$ws.Range("J113").Value = 11033.333
$ws.Range("K113").Value = 3798.5
$ws.Range("L113").Value = 11033.333
$ws.Range("M113").Value = -544.5
$ws.Range("N113").Value = -17541.333
# Row 137
$ws.Range("H137").Value = 1925.7142
$ws.Range("J137").Value = 3666
$ws.Range("L137").Value = 10998
$ws.Range("N137").Value = -16098
# Row 138
$ws.Range("H138").Value = 317642
$ws.Range("J138").Value = 423095.38
$ws.Range("L138").Value = 1269286.14
$ws.Range("N138").Value = -1279566.14

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 6348.077
$ws.Range("I45").Value = 3866
$ws.Range("K45").Value = 3866
$ws.Range("M45").Value = -3489
# Row 62
$ws.Range("H62").Value = 89999
$ws.Range("J62").Value = 89999
$ws.Range("L62").Value = 89999
$ws.Range("N62").Value = -91247
# Row 65
$ws.Range("H65").Value = 89999
$ws.Range("J65").Value = 89999
$ws.Range("L65").Value = 269997
$ws.Range("N65").Value = -276237
# Row 74
$ws.Range("H74").Value = 1502.6086
$ws.Range("I74").Value = 851.75
$ws.Range("K74").Value = 851.75
$ws.Range("M74").Value = 22.25
# Row 75
$ws.Range("H75").Value = 60000
$ws.Range("J75").Value = 60000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61748
# Row 77
$ws.Range("H77").Value = 1502.6086
$ws.Range("I77").Value = 851.75
$ws.Range("K77").Value = 4258.75
$ws.Range("M77").Value = 109.25
# Row 78
$ws.Range("H78").Value = 60000
$ws.Range("J78").Value = 60000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -188736
# Row 132
$ws.Range("H132").Value = 5916.8804
$ws.Range("I132").Value = 4042.1135
$ws.Range("J132").Value = 9503.392
$ws.Range("K132").Value = 12126.3405
$ws.Range("L132").Value = 28510.176
$ws.Range("M132").Value = -9596.3405
$ws.Range("N132").Value = -33570.176

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1559
$ws.Range("I20").Value = 1597.8889
$ws.Range("K20").Value = 1597.8889
$ws.Range("M20").Value = -1350.8889
# Row 25
$ws.Range("H25").Value = 1250.75
$ws.Range("I25").Value = 1250.75
$ws.Range("K25").Value = 1250.75
$ws.Range("M25").Value = -1015.75
# Row 32
$ws.Range("H32").Value = 21000
$ws.Range("J32").Value = 21000
$ws.Range("L32").Value = 21000
$ws.Range("N32").Value = -21768
# Row 105
$ws.Range("H105").Value = 7380.3335
$ws.Range("I105").Value = 6493.2
$ws.Range("K105").Value = 6493.2
$ws.Range("M105").Value = -4746.2

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 227.45454
$ws.Range("I7").Value = 239.16667
$ws.Range("K7").Value = 239.16667
$ws.Range("M7").Value = -126.16667
# Row 31
$ws.Range("H31").Value = 258123.8
$ws.Range("I31").Value = 417871.53
$ws.Range("K31").Value = 417871.53
$ws.Range("M31").Value = -417576.53
# Row 34
$ws.Range("H34").Value = 258123.8
$ws.Range("I34").Value = 417871.53
$ws.Range("K34").Value = 417871.53
$ws.Range("M34").Value = -417669.53
# Row 58
$ws.Range("H58").Value = 2154.16
$ws.Range("I58").Value = 2018.6
$ws.Range("J58").Value = 2357.5
$ws.Range("K58").Value = 2018.6
$ws.Range("L58").Value = 2357.5
$ws.Range("M58").Value = -1815.6
$ws.Range("N58").Value = -2763.5
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null
# Row 99
$ws.Range("H99").Value = 5015.268
$ws.Range("I99").Value = 4811.727
$ws.Range("K99").Value = 4811.727
$ws.Range("M99").Value = -3313.727
# Row 105
$ws.Range("H105").Value = 4138.1113
$ws.Range("I105").Value = 4138.1113
$ws.Range("K105").Value = 4138.1113
$ws.Range("M105").Value = -2391.1113
# Row 126
$ws.Range("H126").Value = 5015.268
$ws.Range("I126").Value = 4811.727
$ws.Range("K126").Value = 14435.181
$ws.Range("M126").Value = -11965.181
# Row 136
$ws.Range("H136").Value = 2154.16
$ws.Range("I136").Value = 2018.6
$ws.Range("J136").Value = 2357.5
$ws.Range("K136").Value = 6055.799999999999
$ws.Range("L136").Value = 7072.5
$ws.Range("M136").Value = -3505.799999999999
$ws.Range("N136").Value = -12172.5

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 333995
$ws.Range("I9").Value = 667333.3
$ws.Range("J9").Value = 656.6667
$ws.Range("K9").Value = 2001999.9
$ws.Range("L9").Value = 1970.0001
$ws.Range("M9").Value = -2001775.9
$ws.Range("N9").Value = -2418.0001
# Row 38
$ws.Range("H38").Value = 2894.1904
$ws.Range("I38").Value = 295
$ws.Range("K38").Value = 885
$ws.Range("M38").Value = -538
# Row 51
$ws.Range("H51").Value = 5601
$ws.Range("I51").Value = 4973.25
$ws.Range("J51").Value = 6103.2
$ws.Range("K51").Value = 14919.75
$ws.Range("L51").Value = 18309.6
$ws.Range("M51").Value = -14459.75
$ws.Range("N51").Value = -19229.6
# Row 92
$ws.Range("H92").Value = 5075
$ws.Range("J92").Value = 5075
$ws.Range("L92").Value = 15225
$ws.Range("N92").Value = -17721
# Row 120
$ws.Range("H120").Value = 7335.8335
$ws.Range("I120").Value = 7335.8335
$ws.Range("K120").Value = 22007.5005
$ws.Range("M120").Value = -17169.5005
# Row 129
$ws.Range("H129").Value = 3569.8096
$ws.Range("J129").Value = 7288.222
$ws.Range("L129").Value = 21864.666
$ws.Range("N129").Value = -31864.666
# Row 131
$ws.Range("H131").Value = 2074.4614
$ws.Range("J131").Value = 2299.2
$ws.Range("L131").Value = 6897.599999999999
$ws.Range("N131").Value = -16977.6

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 36159.793
$ws.Range("I102").Value = 1164.6818
$ws.Range("K102").Value = 1164.6818
$ws.Range("M102").Value = 457.3181999999999
# Row 107
$ws.Range("H107").Value = 1141.3334
$ws.Range("I107").Value = 1419
$ws.Range("J107").Value = 863.6667
$ws.Range("K107").Value = 1419
$ws.Range("L107").Value = 863.6667
$ws.Range("M107").Value = 501
$ws.Range("N107").Value = -4703.6667
# Row 113
$ws.Range("H113").Value = 3750.923
$ws.Range("I113").Value = 4661.5
$ws.Range("J113").Value = 2970.4285
$ws.Range("K113").Value = 4661.5
$ws.Range("L113").Value = 2970.4285
$ws.Range("M113").Value = -2491.5
$ws.Range("N113").Value = -7310.4285

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5316.885
$ws.Range("I7").Value = 4591.5
$ws.Range("J7").Value = 6477.5
$ws.Range("K7").Value = 4591.5
$ws.Range("L7").Value = 6477.5
$ws.Range("M7").Value = -4479.5
$ws.Range("N7").Value = -6701.5
# Row 68
$ws.Range("H68").Value = 2574.4348
$ws.Range("I68").Value = 2440.3333
$ws.Range("K68").Value = 2440.3333
$ws.Range("M68").Value = -1691.3333
# Row 71
$ws.Range("H71").Value = 2574.4348
$ws.Range("I71").Value = 2440.3333
$ws.Range("K71").Value = 12201.6665
$ws.Range("M71").Value = -8457.666499999999
# Row 104
$ws.Range("H104").Value = 48318.89
$ws.Range("J104").Value = 48318.89
$ws.Range("L104").Value = 48318.89
$ws.Range("N104").Value = -55306.89
# Row 126
$ws.Range("H126").Value = 5316.885
$ws.Range("I126").Value = 4591.5
$ws.Range("J126").Value = 6477.5
$ws.Range("K126").Value = 13774.5
$ws.Range("L126").Value = 19432.5
$ws.Range("M126").Value = -11304.5
$ws.Range("N126").Value = -24372.5

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 13192
$ws.Range("J41").Value = 11500.333
$ws.Range("L41").Value = 11500.333
$ws.Range("N41").Value = -12280.333
# Row 122
$ws.Range("H122").Value = 3064.75
$ws.Range("I122").Value = 3451.5
$ws.Range("K122").Value = 10354.5
$ws.Range("M122").Value = -7904.5
# Row 136
$ws.Range("H136").Value = 4506.5674
$ws.Range("I136").Value = 4506.5674
$ws.Range("K136").Value = 13519.7022
$ws.Range("M136").Value = -10969.7022
